$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") for data rows 2-15 changes from serial date 45207 to 45208
for ($r = 2; $r -le 15; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45207) {
        $cell.Value2 = 45208
    }
}
